$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows 344-357 (dates 44418-44431)
# columns: A = date serial, B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$newData = @(
    @(344, 44418, 0, 0, 0),
    @(345, 44419, 0, 0, 0),
    @(346, 44420, 1, 1, 46.70714619336758),
    @(347, 44421, 0, 1, 46.70714619336758),
    @(348, 44422, 0, 1, 46.70714619336758),
    @(349, 44423, 0, 1, 46.70714619336758),
    @(350, 44424, 0, 1, 46.70714619336758),
    @(351, 44425, 0, 1, 46.70714619336758),
    @(352, 44426, 0, 1, 46.70714619336758),
    @(353, 44427, 0, 0, 0),
    @(354, 44428, 0, 0, 0),
    @(355, 44429, 1, 1, 46.70714619336758),
    @(356, 44430, 0, 1, 46.70714619336758),
    @(357, 44431, 0, 1, 46.70714619336758)
)

# Use row 343 as the template for formatting of A, B, C, D columns
$templateRow = 343

foreach ($entry in $newData) {
    $r = $entry[0]
    $dateVal = $entry[1]
    $bVal = $entry[2]
    $cVal = $entry[3]
    $dVal = $entry[4]

    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal

    # Copy formatting (style) from the template row's cells
    $ws.Cells.Item($templateRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($templateRow, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($templateRow, 3).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($templateRow, 4).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
}
